$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C, E contain non-numeric-looking text and can be set directly.
# Column D sometimes contains values that look numeric (e.g. "227.67"); a leading
# apostrophe forces Excel to store them as text, and resetting the Style back to
# "Normal" afterwards strips the quote-prefix formatting Excel applies, so the cell
# keeps its original (default) style while the text value is preserved exactly.

$ws.Range("D2").Value = "'34.431.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "'1.802.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").Value = "'227.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("D6").Value = "'0.580"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.99%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").Value = "'35.01"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.24%  "
$ws.Range("E9").Value = "  +0.49%  "
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").Value = "'2.062.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").Value = "'11.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "'1.807.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").Value = "'34.374.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.43%  "
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "'69.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("D20").Value = "'244.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").Value = "'11.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("D23").Value = "'4.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").Value = "'170.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.58%  "
$ws.Range("D25").Value = "'2.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.11%  "
$ws.Range("E26").Value = "  +3.97%  "
$ws.Range("D27").Value = "'16.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("D28").Value = "'0.118"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "'3.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("D33").Value = "'3.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.23%  "
$ws.Range("D34").Value = "'1.82"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "'1.398.18"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.85%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.679"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.85%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'2.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("D40").Value = "'83.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("E41").Value = "  +3.38%  "
$ws.Range("D42").Value = "'0.946"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("D43").Value = "'2.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.28%  "
$ws.Range("D44").Value = "'13.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.15%  "
$ws.Range("D45").Value = "'1.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("E46").Value = "  -2.34%  "
$ws.Range("D47").Value = "'5.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.42%  "
$ws.Range("D48").Value = "'1.964.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").Value = "'104.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.72%  "
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("E51").Value = "  -0.43%  "
